# "validacao de dados, avisos e informacoes"
#
# Adds a new worksheet "Aula2" (a copy of the "Base" sheet layout), placed
# right after "Aula 1", pre-filled with team values in column B and a
# teams list in column H (swapping "Mega" for "Black" in H3), and applies
# a list-type data validation on B2:B15 driven by the H2:H5 list.

$wb = $excel.ActiveWorkbook

# Duplicate the "Base" sheet (same header/row/style skeleton is reused by
# the new "Aula2" sheet) and drop the copy right after "Aula 1".
$base = $wb.Worksheets.Item("Base")
$aula1 = $wb.Worksheets.Item("Aula 1")
$base.Copy($null, $aula1)

$ws = $wb.Worksheets.Item(2)
$ws.Name = "Aula2"

# H2:H4 already carry Ultra/Mega/Blaster/Power (copied from "Base"); only
# H3 needs to change from "Mega" to the new team "Black".
$ws.Range("H3").Value = "Black"

# Fill in the team assigned to each of the 14 students in column B.
$teams = @("Ultra", "Blaster", "Blaster", "Mega", "Power", "Mega", "Ultra", "Blaster", "Power", "Ultra", "Power", "Mega", "Ultra", "Blaster")
for ($i = 0; $i -lt $teams.Length; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $teams[$i]
}

# Dropdown list validation on the team column, sourced from H2:H5.
$ws.Range("B2:B15").Validation.Add(3, 1, 1, "=`$H`$2:`$H`$5")

$ws.Range("B22").Select()
